$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading-percent values for rows 2-25 (A = 0..23), columns B:F and J:M.
# Data laid out as one row per array element: B,C,D,E,F,J,K,L,M
$data = @(
    @(21.00161372784055,2.971355352184129,4.486033468703653,10.53686874755594,71.63293216806524,11.31345065817434,17.29396089744464,10.81972464555603,18.70543065278988),
    @(21.03016882288832,2.819988934432768,4.356415891133224,10.55991303917255,70.82504543664436,11.30218204383513,17.33297963058629,10.85970154143504,18.76925023135885),
    @(21.05400610868024,2.776838938734645,4.273843329462036,10.57490195254308,70.32644243018918,11.2953799438293,17.36201502296626,10.88579955826696,18.81208635761337),
    @(21.06530247405438,2.765640278106161,4.239464926892093,10.58122186789461,70.12274315558268,11.29263745089993,17.3751208662343,10.89682589449631,18.8304601804316),
    @(21.067273707782,2.763815866735554,4.233713036038391,10.58228409795148,70.08889171700487,11.29218383947566,17.37737390849612,10.8986804625368,18.83356656266567),
    @(21.05415205235955,2.776685570714607,4.273382613565515,10.57498632660982,70.32369718454717,11.29534283845374,17.3621866203903,10.88594667821422,18.81233043782434),
    @(21.0101496012547,2.920193009801353,4.441973625300726,10.54464056480838,71.35498077092858,11.30954099955497,17.30635890308263,10.83318709888852,18.72667781879594),
    @(20.97397554205978,3.27019218309845,4.748030900005522,10.49176393706377,73.35140493925365,11.33831069465592,17.2372936729259,10.74199967621628,18.58768980651897),
    @(20.97802747226157,3.502869419751057,4.956923533248705,10.45691557526518,74.79490540097549,11.36001433475531,17.21132990371408,10.68243069661502,18.50325401577514),
    @(20.98652062791047,3.60336083748751,5.048310776973781,10.44192171643683,75.44497729968231,11.37001109650154,17.20491903955034,10.65693205546888,18.46868341311757),
    @(20.99069138652743,3.640642027955948,5.082380787547298,10.43636673720712,75.6900669662371,11.37381433758496,17.20326872315728,10.64750551947999,18.45614477409003),
    @(20.98975071984342,3.63264725557303,5.0750672477273,10.43755764669207,75.63733255273111,11.37299445674557,17.20358956905592,10.64952550904312,18.45882061645348),
    @(20.98684464021928,3.60644348178396,5.051124556063568,10.44146224605163,75.4651633322563,11.37032363761738,17.204767685371,10.65615193862633,18.46764077691112),
    @(20.98518883264373,3.590292219692845,5.036388754943149,10.44386990855179,75.35956040581006,11.3686899836928,17.20559056246333,10.66024064997909,18.47311534153991),
    @(20.97760611640891,3.496194059895034,4.95087678327253,10.45791268722469,74.75227786946884,11.35936350316809,17.21185761499082,10.6841292141987,18.50559058117039),
    @(20.97465667352009,3.437094115652996,4.897476125670179,10.46674699587579,74.37794887930625,11.35367361675759,17.21708610615924,10.69919321303744,18.52649678337206),
    @(20.97358631086865,3.402597471053245,4.866419634195828,10.47190913161473,74.16203277682611,11.35041259132642,17.22060167348859,10.70800823277224,18.53888286557605),
    @(20.9733314697793,3.390831053862734,4.855846150616131,10.47367085170838,74.08882624668975,11.34931046732119,17.22187924802564,10.71101873746299,18.5431386373789),
    @(20.97490585928081,3.443437582778115,4.903196195891471,10.46579820296736,74.41786113888169,11.35427811014513,17.21647691218832,10.69757404318007,18.52423387762429),
    @(20.98767233688898,3.614161153217334,5.058171758139559,10.44031204091876,75.51576391031711,11.37110764283071,17.20440054364095,10.65419938004844,18.46503508463308),
    @(21.00157860033015,3.72123430996896,5.156325941460818,10.4243712504129,76.22696178893798,11.38220973984455,17.20103884649515,10.62718740266648,18.4295657766073),
    @(20.993648371027,3.664500141267585,5.104229651681388,10.43281385322377,75.84800543010253,11.37627495043754,17.20241834680772,10.64148222265872,18.44820162479974),
    @(20.9747912544948,3.440571318622836,4.900611259250849,10.46622689311919,74.39981900576366,11.35400478702384,17.21675074136763,10.69830558847572,18.52525579469679),
    @(20.97838506227171,3.179783786733471,4.667972888190977,10.50536295188881,72.81486936727242,11.33042828990342,17.25163462422715,10.76536012803996,18.6221860916994)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E
    $ws.Cells.Item($row, 6).Value = $vals[4]   # F
    $ws.Cells.Item($row, 10).Value = $vals[5]  # J
    $ws.Cells.Item($row, 11).Value = $vals[6]  # K
    $ws.Cells.Item($row, 12).Value = $vals[7]  # L
    $ws.Cells.Item($row, 13).Value = $vals[8]  # M
}
